{"js": "// Replace each old equation text with its updated equation text.\n// Mapping was derived 1:1 (each old value is unique in the document).\nconst replacements = [\n  [\"38-21=17\", \"20+56=76\"],\n  [\"1+30=31\", \"48+2=50\"],\n  [\"12+23=35\", \"35-7=28\"],\n  [\"52-34=18\", \"66-44=22\"],\n  [\"6+67=73\", \"30+67=97\"],\n  [\"26+41=67\", \"88-48=40\"],\n  [\"99-6=93\", \"79-70=9\"],\n  [\"91-47=44\", \"0+53=53\"],\n  [\"12+38=50\", \"30+29=59\"],\n  [\"58+32=90\", \"88-6=82\"],\n  [\"16+12=28\", \"59-24=35\"],\n  [\"29+12=41\", \"7+46=53\"],\n  [\"44-8=36\", \"52+34=86\"],\n  [\"96-2=94\", \"10+26=36\"],\n  [\"2+87=89\", \"32-12=20\"],\n  [\"72+27=99\", \"95-58=37\"],\n  [\"43+53=96\", \"20-3=17\"],\n  [\"87-47=40\", \"44+27=71\"],\n  [\"6-0=6\", \"44+17=61\"],\n  [\"7+62=69\", \"51+47=98\"],\n  [\"25+3=28\", \"76-74=2\"],\n  [\"94-64=30\", \"87-74=13\"],\n  [\"13-4=9\", \"80-15=65\"],\n  [\"6+22=28\", \"19+65=84\"],\n  [\"81+11=92\", \"9+81=90\"],\n  [\"10+2=12\", \"70+2=72\"],\n  [\"38+49=87\", \"4+70=74\"],\n  [\"70+0=70\", \"79+16=95\"],\n  [\"9+71=80\", \"55-42=13\"],\n  [\"19-7=12\", \"71+21=92\"],\n  [\"7+25=32\", \"82-18=64\"],\n  [\"82-1=81\", \"23-7=16\"],\n  [\"47+30=77\", \"94+5=99\"],\n  [\"33-13=20\", \"60-31=29\"],\n  [\"19+49=68\", \"31+29=60\"],\n  [\"84-29=55\", \"57+13=70\"],\n  [\"73-63=10\", \"12+24=36\"],\n  [\"34+13=47\", \"15-5=10\"],\n  [\"56+40=96\", \"4+78=82\"],\n  [\"40+52=92\", \"92-18=74\"],\n  [\"99-69=30\", \"59-9=50\"],\n  [\"64-4=60\", \"24+15=39\"],\n  [\"86-59=27\", \"23-22=1\"],\n  [\"56+42=98\", \"48+9=57\"],\n  [\"51+21=72\", \"44-26=18\"],\n  [\"69-16=53\", \"84+4=88\"],\n  [\"72+6=78\", \"14+67=81\"],\n  [\"1+95=96\", \"3+71=74\"],\n  [\"7+38=45\", \"93-81=12\"],\n  [\"95-67=28\", \"86-21=65\"],\n  [\"14+12=26\", \"20+17=37\"],\n  [\"37+55=92\", \"23+58=81\"],\n  [\"28-24=4\", \"96-0=96\"],\n  [\"48+20=68\", \"74-16=58\"],\n  [\"16+18=34\", \"93-53=40\"],\n  [\"23+15=38\", \"5+23=28\"],\n  [\"50+4=54\", \"58+9=67\"],\n  [\"85+12=97\", \"77-1=76\"],\n  [\"86-54=32\", \"88-51=37\"],\n  [\"47+34=81\", \"90-14=76\"],\n  [\"62+12=74\", \"41+13=54\"],\n  [\"33+52=85\", \"82+1=83\"],\n  [\"49-45=4\", \"7+17=24\"],\n  [\"83+3=86\", \"59+23=82\"],\n  [\"7+37=44\", \"73-14=59\"],\n  [\"80-78=2\", \"40+11=51\"],\n  [\"90-21=69\", \"17+31=48\"],\n  [\"35+37=72\", \"99-12=87\"],\n  [\"58-41=17\", \"49+35=84\"],\n  [\"0+94=94\", \"71-18=53\"],\n  [\"41-35=6\", \"68-47=21\"],\n  [\"67-39=28\", \"19+61=80\"],\n  [\"10+22=32\", \"4+81=85\"],\n  [\"22+34=56\", \"76+9=85\"],\n  [\"72-60=12\", \"97-75=22\"],\n  [\"21+30=51\", \"29+29=58\"],\n  [\"62-38=24\", \"12+19=31\"],\n  [\"42-33=9\", \"34+42=76\"],\n  [\"90+8=98\", \"44-28=16\"],\n  [\"11-4=7\", \"9+0=9\"],\n  [\"77-61=16\", \"96-83=13\"],\n  [\"48-37=11\", \"20+22=42\"],\n  [\"26+8=34\", \"55+16=71\"],\n  [\"99-36=63\", \"27+70=97\"],\n  [\"90-64=26\", \"0+27=27\"],\n  [\"97-2=95\", \"75-32=43\"],\n  [\"96-39=57\", \"10+6=16\"],\n  [\"54+6=60\", \"93-7=86\"],\n  [\"50-48=2\", \"27+18=45\"],\n  [\"39+6=45\", \"88-8=80\"],\n  [\"57-39=18\", \"12+27=39\"],\n  [\"13+20=33\", \"94-50=44\"],\n  [\"37+24=61\", \"28+6=34\"],\n  [\"63-55=8\", \"86+4=90\"],\n  [\"6+84=90\", \"81-40=41\"],\n  [\"45+43=88\", \"23+35=58\"],\n  [\"56-8=48\", \"32+27=59\"],\n  [\"96-11=85\", \"95-91=4\"],\n  [\"51-0=51\", \"1+93=94\"],\n  [\"52+29=81\", \"92-83=9\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error('No match found for: ' + oldText);\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();", "ps1": "# Replace each old equation text with its updated equation text.\n# Mapping was derived 1:1 (each old value is unique in the document).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"38-21=17\", \"20+56=76\")\n    ,@(\"1+30=31\", \"48+2=50\")\n    ,@(\"12+23=35\", \"35-7=28\")\n    ,@(\"52-34=18\", \"66-44=22\")\n    ,@(\"6+67=73\", \"30+67=97\")\n    ,@(\"26+41=67\", \"88-48=40\")\n    ,@(\"99-6=93\", \"79-70=9\")\n    ,@(\"91-47=44\", \"0+53=53\")\n    ,@(\"12+38=50\", \"30+29=59\")\n    ,@(\"58+32=90\", \"88-6=82\")\n    ,@(\"16+12=28\", \"59-24=35\")\n    ,@(\"29+12=41\", \"7+46=53\")\n    ,@(\"44-8=36\", \"52+34=86\")\n    ,@(\"96-2=94\", \"10+26=36\")\n    ,@(\"2+87=89\", \"32-12=20\")\n    ,@(\"72+27=99\", \"95-58=37\")\n    ,@(\"43+53=96\", \"20-3=17\")\n    ,@(\"87-47=40\", \"44+27=71\")\n    ,@(\"6-0=6\", \"44+17=61\")\n    ,@(\"7+62=69\", \"51+47=98\")\n    ,@(\"25+3=28\", \"76-74=2\")\n    ,@(\"94-64=30\", \"87-74=13\")\n    ,@(\"13-4=9\", \"80-15=65\")\n    ,@(\"6+22=28\", \"19+65=84\")\n    ,@(\"81+11=92\", \"9+81=90\")\n    ,@(\"10+2=12\", \"70+2=72\")\n    ,@(\"38+49=87\", \"4+70=74\")\n    ,@(\"70+0=70\", \"79+16=95\")\n    ,@(\"9+71=80\", \"55-42=13\")\n    ,@(\"19-7=12\", \"71+21=92\")\n    ,@(\"7+25=32\", \"82-18=64\")\n    ,@(\"82-1=81\", \"23-7=16\")\n    ,@(\"47+30=77\", \"94+5=99\")\n    ,@(\"33-13=20\", \"60-31=29\")\n    ,@(\"19+49=68\", \"31+29=60\")\n    ,@(\"84-29=55\", \"57+13=70\")\n    ,@(\"73-63=10\", \"12+24=36\")\n    ,@(\"34+13=47\", \"15-5=10\")\n    ,@(\"56+40=96\", \"4+78=82\")\n    ,@(\"40+52=92\", \"92-18=74\")\n    ,@(\"99-69=30\", \"59-9=50\")\n    ,@(\"64-4=60\", \"24+15=39\")\n    ,@(\"86-59=27\", \"23-22=1\")\n    ,@(\"56+42=98\", \"48+9=57\")\n    ,@(\"51+21=72\", \"44-26=18\")\n    ,@(\"69-16=53\", \"84+4=88\")\n    ,@(\"72+6=78\", \"14+67=81\")\n    ,@(\"1+95=96\", \"3+71=74\")\n    ,@(\"7+38=45\", \"93-81=12\")\n    ,@(\"95-67=28\", \"86-21=65\")\n    ,@(\"14+12=26\", \"20+17=37\")\n    ,@(\"37+55=92\", \"23+58=81\")\n    ,@(\"28-24=4\", \"96-0=96\")\n    ,@(\"48+20=68\", \"74-16=58\")\n    ,@(\"16+18=34\", \"93-53=40\")\n    ,@(\"23+15=38\", \"5+23=28\")\n    ,@(\"50+4=54\", \"58+9=67\")\n    ,@(\"85+12=97\", \"77-1=76\")\n    ,@(\"86-54=32\", \"88-51=37\")\n    ,@(\"47+34=81\", \"90-14=76\")\n    ,@(\"62+12=74\", \"41+13=54\")\n    ,@(\"33+52=85\", \"82+1=83\")\n    ,@(\"49-45=4\", \"7+17=24\")\n    ,@(\"83+3=86\", \"59+23=82\")\n    ,@(\"7+37=44\", \"73-14=59\")\n    ,@(\"80-78=2\", \"40+11=51\")\n    ,@(\"90-21=69\", \"17+31=48\")\n    ,@(\"35+37=72\", \"99-12=87\")\n    ,@(\"58-41=17\", \"49+35=84\")\n    ,@(\"0+94=94\", \"71-18=53\")\n    ,@(\"41-35=6\", \"68-47=21\")\n    ,@(\"67-39=28\", \"19+61=80\")\n    ,@(\"10+22=32\", \"4+81=85\")\n    ,@(\"22+34=56\", \"76+9=85\")\n    ,@(\"72-60=12\", \"97-75=22\")\n    ,@(\"21+30=51\", \"29+29=58\")\n    ,@(\"62-38=24\", \"12+19=31\")\n    ,@(\"42-33=9\", \"34+42=76\")\n    ,@(\"90+8=98\", \"44-28=16\")\n    ,@(\"11-4=7\", \"9+0=9\")\n    ,@(\"77-61=16\", \"96-83=13\")\n    ,@(\"48-37=11\", \"20+22=42\")\n    ,@(\"26+8=34\", \"55+16=71\")\n    ,@(\"99-36=63\", \"27+70=97\")\n    ,@(\"90-64=26\", \"0+27=27\")\n    ,@(\"97-2=95\", \"75-32=43\")\n    ,@(\"96-39=57\", \"10+6=16\")\n    ,@(\"54+6=60\", \"93-7=86\")\n    ,@(\"50-48=2\", \"27+18=45\")\n    ,@(\"39+6=45\", \"88-8=80\")\n    ,@(\"57-39=18\", \"12+27=39\")\n    ,@(\"13+20=33\", \"94-50=44\")\n    ,@(\"37+24=61\", \"28+6=34\")\n    ,@(\"63-55=8\", \"86+4=90\")\n    ,@(\"6+84=90\", \"81-40=41\")\n    ,@(\"45+43=88\", \"23+35=58\")\n    ,@(\"56-8=48\", \"32+27=59\")\n    ,@(\"96-11=85\", \"95-91=4\")\n    ,@(\"51-0=51\", \"1+93=94\")\n    ,@(\"52+29=81\", \"92-83=9\")\n)\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceAll)\n    if (-not $found) {\n        Write-Output \"WARNING: no match for $oldText\"\n    }\n}"}
